# Append a new weekly schedule entry (period 84) as row 34 on 工作表1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = "2025/11/14"
$ws.Range("B34").Value = "2026/1/9"
$ws.Range("C34").Value = "第84期 秘寶 開放區域 殤金國 祕寶效果: 核心進階傷害提高1534930 (11051496)"

# Mirror the cursor position left behind after typing the new row in Excel.
$ws.Range("C37").Select()
